$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "description" column right after "name" (column A), pushing
# the existing columns (parent_sku, sku, price, ...) one to the right.
$ws.Columns("B").Insert()
$ws.Range("B1").Value = "description"
$ws.Range("B2").Value = "Product description (Can be in HTML)"
$ws.Range("B3").Value = "Product description (Can be in HTML)"

# The "new" / "active" / "manage_stock" flags are now stored as plain
# numbers instead of booleans.
$ws.Range("K2").Value = 1
$ws.Range("M2").Value = 1
$ws.Range("N2").Value = 1
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 1

# The column insert shifts the "images" hyperlink data to column L, but
# leaves the hyperlink anchor behind at K - rebuild the hyperlinks so they
# point at the cell that now actually holds the URL text.
$ws.Range("K2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("L2"), "http://kommercio.id/assets/images/logo.png;http:/kommercio.id/assets/images/logo.png")
$ws.Hyperlinks.Add($ws.Range("L3"), "http://kommercio.id/assets/images/logo.png;http:/kommercio.id/assets/images/logo.png")

# Append a new "created_at" column at the end of the table.
$ws.Range("O1").Copy()
$ws.Range("P1").PasteSpecial(-4122)
$ws.Range("P1").Value = "created_at"
$ws.Range("P2").Value = "2017-05-15 10:00:00"
$ws.Range("P3").Value = "2017-05-15 10:00:00"

# Selection now covers a single cell instead of B2:B3.
$ws.Range("B2").Select()
